## edit.ps1 - applies the "List of Tasks - Notes.docx" revision described by
## the supplied XML diff.
##
## Strategy: Word's plain Find/Replace collapses everything into a single
## run, which does not match the target (which has runs split around
## spell-checked words with w:proofErr markers, merged runs where multiple
## runs collapse into one, and runs removed entirely). To get byte-accurate
## run structure we rebuild the *entire* run content of each paragraph that
## changes, and push it in with Range.InsertXML (a WordprocessingML package
## fragment) over the paragraph's full content range (excluding the
## paragraph mark). This was verified empirically to be the only reliable
## way to get exact run/proofErr structure with this host - partial
## mid-paragraph InsertXML calls land content at the wrong offset.

$d = $word.ActiveDocument

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphRuns($Index, $InnerXml) {
    $p = $d.Paragraphs($Index)
    $start = $p.Range.Start
    $end = $p.Range.End
    $rng = $d.Range($start, $end - 1)
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wordNs + '><w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

## 1. "Research other HardwareX Papers" (strike) -> split around "HardwareX"
##    with proofErr spell-check markers.
$xml1 = (
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Research other </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t>HardwareX</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> Papers</w:t></w:r>'
)
Set-ParagraphRuns 2 $xml1

## 2. "Add the HardwareX LaTeX Files" (strike) -> split around "HardwareX".
$xml2 = (
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Add the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t>HardwareX</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> LaTeX Files</w:t></w:r>'
)
Set-ParagraphRuns 13 $xml2

## 3. "Setup Zenodo with eDNA GitHub (Maybe on the UI and Framework GitHub
##    as well)" -> split around "Zenodo".
$xml3 = (
    '<w:r><w:t xml:space="preserve">Setup </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Zenodo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> with eDNA GitHub (Maybe on the UI and Framework GitHub as well)</w:t></w:r>'
)
Set-ParagraphRuns 17 $xml3

## 4. "Hardware in Context" paragraph: merge the three runs " Draft
##    Completed" + ", " + "need to prep for Godshalk" into a single run.
$xml4 = (
    '<w:r><w:t xml:space="preserve">Hardware in Context – Kai </w:t></w:r>' +
    '<w:r><w:t>(</w:t></w:r>' +
    '<w:r><w:t>1</w:t></w:r>' +
    '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Draft Completed, need to prep for Godshalk</w:t></w:r>' +
    '<w:r><w:t>)</w:t></w:r>'
)
Set-ParagraphRuns 23 $xml4

## 5. "Software Setup Instructions - Nathan/Jacob/Mark/Kai" ->
##    "Software Setup Instructions - Mark".
$xml5 = (
    '<w:r><w:t>Software Setup Instructions –</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>Mark</w:t></w:r>'
)
Set-ParagraphRuns 27 $xml5

## 6. "Validation and Characterization - Riley/Hendy/Kai" -> drop the "/Kai"
##    run.
$xml6 = '<w:r><w:t>Validation and Characterization – Riley/Hendy</w:t></w:r>'
Set-ParagraphRuns 28 $xml6

## 7. "Conclusion - Riley/Kai" -> "Conclusion - Riley".
$xml7 = '<w:r><w:t>Conclusion – Riley</w:t></w:r>'
Set-ParagraphRuns 29 $xml7

Write-Host "Done."
